$wb = $excel.ActiveWorkbook

# --- Step 1: insert a new sheet "2022-Q1" (copy of 2021-Q1 structure) right before the "总计" sheet ---
# NB: sheet variables track by position, not stable identity, in this host -- so we
# capture the insertion index, perform the copy, then re-resolve handles by that index/name.
$insertIdx = $wb.Worksheets.Count
$totalSheetRef = $wb.Worksheets.Item($insertIdx)
$template = $wb.Worksheets.Item("2021-Q1")
$template.Copy($totalSheetRef)
$newSheet = $wb.Worksheets.Item($insertIdx)
$newSheet.Name = "2022-Q1"

# Fix the D1 header text (template sheets use "基金金额"; new sheet uses "基金规模")
$newSheet.Cells.Item(1,4).Value = "基金规模"

# --- Step 2: write the fund holding detail rows into "2022-Q1" ---
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = '''008985'
$newSheet.Cells.Item(2,3).Value = '东方红启东三年持有期混合'
$newSheet.Cells.Item(2,4).Value = '''130.09'
$newSheet.Cells.Item(2,5).Value = '''90.89'
$newSheet.Cells.Item(2,6).Value = '''3.89'
$newSheet.Cells.Item(2,7).Value = '''5.0605'
$newSheet.Cells.Item(2,8).Value = 5
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = '''010902'
$newSheet.Cells.Item(3,3).Value = '博时成长领航灵活配置混合A'
$newSheet.Cells.Item(3,4).Value = '''62.64'
$newSheet.Cells.Item(3,5).Value = '''73.89'
$newSheet.Cells.Item(3,6).Value = '''5.64'
$newSheet.Cells.Item(3,7).Value = '''3.5329'
$newSheet.Cells.Item(3,8).Value = 2
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = '''012344'
$newSheet.Cells.Item(4,3).Value = '嘉实领先优势混合型证券投资基金A'
$newSheet.Cells.Item(4,4).Value = '''68.76'
$newSheet.Cells.Item(4,5).Value = '''82.99'
$newSheet.Cells.Item(4,6).Value = '''4.70'
$newSheet.Cells.Item(4,7).Value = '''3.2317'
$newSheet.Cells.Item(4,8).Value = 5
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = '''007802'
$newSheet.Cells.Item(5,3).Value = '兴全合泰混合A'
$newSheet.Cells.Item(5,4).Value = '''83.41'
$newSheet.Cells.Item(5,5).Value = '''90.95'
$newSheet.Cells.Item(5,6).Value = '''3.17'
$newSheet.Cells.Item(5,7).Value = '''2.6441'
$newSheet.Cells.Item(5,8).Value = 9
$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = '''010041'
$newSheet.Cells.Item(6,3).Value = '嘉实港股优势混合A'
$newSheet.Cells.Item(6,4).Value = '''49.98'
$newSheet.Cells.Item(6,5).Value = '''92.41'
$newSheet.Cells.Item(6,6).Value = '''4.82'
$newSheet.Cells.Item(6,7).Value = '''2.4090'
$newSheet.Cells.Item(6,8).Value = 6
$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = '''160527'
$newSheet.Cells.Item(7,3).Value = '博时研究优选3年封闭运作灵活配置混合A'
$newSheet.Cells.Item(7,4).Value = '''20.12'
$newSheet.Cells.Item(7,5).Value = '''95.51'
$newSheet.Cells.Item(7,6).Value = '''9.37'
$newSheet.Cells.Item(7,7).Value = '''1.8852'
$newSheet.Cells.Item(7,8).Value = 2
$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = '''014639'
$newSheet.Cells.Item(8,3).Value = '兴证全球合衡三年持有混合A'
$newSheet.Cells.Item(8,4).Value = '''55.51'
$newSheet.Cells.Item(8,5).Value = '''68.89'
$newSheet.Cells.Item(8,6).Value = '''2.25'
$newSheet.Cells.Item(8,7).Value = '''1.2490'
$newSheet.Cells.Item(8,8).Value = 8
$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = '''001878'
$newSheet.Cells.Item(9,3).Value = '嘉实沪港深精选股票'
$newSheet.Cells.Item(9,4).Value = '''23.17'
$newSheet.Cells.Item(9,5).Value = '''93.29'
$newSheet.Cells.Item(9,6).Value = '''5.15'
$newSheet.Cells.Item(9,7).Value = '''1.1933'
$newSheet.Cells.Item(9,8).Value = 5
$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = '''009138'
$newSheet.Cells.Item(10,3).Value = '嘉实瑞成两年持有期混合A'
$newSheet.Cells.Item(10,4).Value = '''22.27'
$newSheet.Cells.Item(10,5).Value = '''75.95'
$newSheet.Cells.Item(10,6).Value = '''4.62'
$newSheet.Cells.Item(10,7).Value = '''1.0289'
$newSheet.Cells.Item(10,8).Value = 2
$newSheet.Cells.Item(11,1).Value = 9
$newSheet.Cells.Item(11,2).Value = '''160726'
$newSheet.Cells.Item(11,3).Value = '嘉实瑞享定期开放灵活配置混合'
$newSheet.Cells.Item(11,4).Value = '''23.58'
$newSheet.Cells.Item(11,5).Value = '''63.95'
$newSheet.Cells.Item(11,6).Value = '''3.98'
$newSheet.Cells.Item(11,7).Value = '''0.9385'
$newSheet.Cells.Item(11,8).Value = 3
$newSheet.Cells.Item(12,1).Value = 10
$newSheet.Cells.Item(12,2).Value = '''910022'
$newSheet.Cells.Item(12,3).Value = '东方红启航三年持有期混合A'
$newSheet.Cells.Item(12,4).Value = '''22.15'
$newSheet.Cells.Item(12,5).Value = '''92.82'
$newSheet.Cells.Item(12,6).Value = '''3.91'
$newSheet.Cells.Item(12,7).Value = '''0.8661'
$newSheet.Cells.Item(12,8).Value = 5
$newSheet.Cells.Item(13,1).Value = 11
$newSheet.Cells.Item(13,2).Value = '''012463'
$newSheet.Cells.Item(13,3).Value = '博时成长优势混合型证券投资基金A'
$newSheet.Cells.Item(13,4).Value = '''19.57'
$newSheet.Cells.Item(13,5).Value = '''75.80'
$newSheet.Cells.Item(13,6).Value = '''3.55'
$newSheet.Cells.Item(13,7).Value = '''0.6947'
$newSheet.Cells.Item(13,8).Value = 9
$newSheet.Cells.Item(14,1).Value = 12
$newSheet.Cells.Item(14,2).Value = '''008966'
$newSheet.Cells.Item(14,3).Value = '博时成长优选两年封闭运作灵活配置混合A'
$newSheet.Cells.Item(14,4).Value = '''7.62'
$newSheet.Cells.Item(14,5).Value = '''80.79'
$newSheet.Cells.Item(14,6).Value = '''8.96'
$newSheet.Cells.Item(14,7).Value = '''0.6828'
$newSheet.Cells.Item(14,8).Value = 1
$newSheet.Cells.Item(15,1).Value = 13
$newSheet.Cells.Item(15,2).Value = '''910028'
$newSheet.Cells.Item(15,3).Value = '东方红内需增长混合型证券投资基金A'
$newSheet.Cells.Item(15,4).Value = '''13.44'
$newSheet.Cells.Item(15,5).Value = '''92.82'
$newSheet.Cells.Item(15,6).Value = '''3.91'
$newSheet.Cells.Item(15,7).Value = '''0.5255'
$newSheet.Cells.Item(15,8).Value = 5
$newSheet.Cells.Item(16,1).Value = 14
$newSheet.Cells.Item(16,2).Value = '''010225'
$newSheet.Cells.Item(16,3).Value = '东方红启航三年持有期混合B'
$newSheet.Cells.Item(16,4).Value = '''13.27'
$newSheet.Cells.Item(16,5).Value = '''92.82'
$newSheet.Cells.Item(16,6).Value = '''3.91'
$newSheet.Cells.Item(16,7).Value = '''0.5189'
$newSheet.Cells.Item(16,8).Value = 5
$newSheet.Cells.Item(17,1).Value = 15
$newSheet.Cells.Item(17,2).Value = '''010903'
$newSheet.Cells.Item(17,3).Value = '博时成长领航灵活配置混合C'
$newSheet.Cells.Item(17,4).Value = '''8.01'
$newSheet.Cells.Item(17,5).Value = '''73.89'
$newSheet.Cells.Item(17,6).Value = '''5.64'
$newSheet.Cells.Item(17,7).Value = '''0.4518'
$newSheet.Cells.Item(17,8).Value = 2
$newSheet.Cells.Item(18,1).Value = 16
$newSheet.Cells.Item(18,2).Value = '''009591'
$newSheet.Cells.Item(18,3).Value = '博时研究精选一年持有期灵活配置混合A'
$newSheet.Cells.Item(18,4).Value = '''8.55'
$newSheet.Cells.Item(18,5).Value = '''76.70'
$newSheet.Cells.Item(18,6).Value = '''4.89'
$newSheet.Cells.Item(18,7).Value = '''0.4181'
$newSheet.Cells.Item(18,8).Value = 4
$newSheet.Cells.Item(19,1).Value = 17
$newSheet.Cells.Item(19,2).Value = '''007803'
$newSheet.Cells.Item(19,3).Value = '兴全合泰混合C'
$newSheet.Cells.Item(19,4).Value = '''12.68'
$newSheet.Cells.Item(19,5).Value = '''90.95'
$newSheet.Cells.Item(19,6).Value = '''3.17'
$newSheet.Cells.Item(19,7).Value = '''0.4020'
$newSheet.Cells.Item(19,8).Value = 9
$newSheet.Cells.Item(20,1).Value = 18
$newSheet.Cells.Item(20,2).Value = '''011740'
$newSheet.Cells.Item(20,3).Value = '博时成长精选混合A'
$newSheet.Cells.Item(20,4).Value = '''7.32'
$newSheet.Cells.Item(20,5).Value = '''75.84'
$newSheet.Cells.Item(20,6).Value = '''5.38'
$newSheet.Cells.Item(20,7).Value = '''0.3938'
$newSheet.Cells.Item(20,8).Value = 4
$newSheet.Cells.Item(21,1).Value = 19
$newSheet.Cells.Item(21,2).Value = '''010042'
$newSheet.Cells.Item(21,3).Value = '嘉实港股优势混合C'
$newSheet.Cells.Item(21,4).Value = '''5.78'
$newSheet.Cells.Item(21,5).Value = '''92.41'
$newSheet.Cells.Item(21,6).Value = '''4.82'
$newSheet.Cells.Item(21,7).Value = '''0.2786'
$newSheet.Cells.Item(21,8).Value = 6
$newSheet.Cells.Item(22,1).Value = 20
$newSheet.Cells.Item(22,2).Value = '''013123'
$newSheet.Cells.Item(22,3).Value = '汇添富精选核心优势一年持有混合A'
$newSheet.Cells.Item(22,4).Value = '''6.15'
$newSheet.Cells.Item(22,5).Value = '''66.61'
$newSheet.Cells.Item(22,6).Value = '''4.26'
$newSheet.Cells.Item(22,7).Value = '''0.2620'
$newSheet.Cells.Item(22,8).Value = 4
$newSheet.Cells.Item(23,1).Value = 21
$newSheet.Cells.Item(23,2).Value = '''002653'
$newSheet.Cells.Item(23,3).Value = '泰康沪港深精选灵活配置混合'
$newSheet.Cells.Item(23,4).Value = '''7.89'
$newSheet.Cells.Item(23,5).Value = '''87.05'
$newSheet.Cells.Item(23,6).Value = '''2.58'
$newSheet.Cells.Item(23,7).Value = '''0.2036'
$newSheet.Cells.Item(23,8).Value = 5
$newSheet.Cells.Item(24,1).Value = 22
$newSheet.Cells.Item(24,2).Value = '''009139'
$newSheet.Cells.Item(24,3).Value = '嘉实瑞成两年持有期混合C'
$newSheet.Cells.Item(24,4).Value = '''4.37'
$newSheet.Cells.Item(24,5).Value = '''75.95'
$newSheet.Cells.Item(24,6).Value = '''4.62'
$newSheet.Cells.Item(24,7).Value = '''0.2019'
$newSheet.Cells.Item(24,8).Value = 2
$newSheet.Cells.Item(25,1).Value = 23
$newSheet.Cells.Item(25,2).Value = '''005228'
$newSheet.Cells.Item(25,3).Value = '汇添富港股通专注成长混合'
$newSheet.Cells.Item(25,4).Value = '''3.64'
$newSheet.Cells.Item(25,5).Value = '''80.11'
$newSheet.Cells.Item(25,6).Value = '''5.21'
$newSheet.Cells.Item(25,7).Value = '''0.1896'
$newSheet.Cells.Item(25,8).Value = 2
$newSheet.Cells.Item(26,1).Value = 24
$newSheet.Cells.Item(26,2).Value = '''012243'
$newSheet.Cells.Item(26,3).Value = '东方红内需增长混合型证券投资基金C'
$newSheet.Cells.Item(26,4).Value = '''4.61'
$newSheet.Cells.Item(26,5).Value = '''92.82'
$newSheet.Cells.Item(26,6).Value = '''3.91'
$newSheet.Cells.Item(26,7).Value = '''0.1803'
$newSheet.Cells.Item(26,8).Value = 5
$newSheet.Cells.Item(27,1).Value = 25
$newSheet.Cells.Item(27,2).Value = '''014640'
$newSheet.Cells.Item(27,3).Value = '兴证全球合衡三年持有混合C'
$newSheet.Cells.Item(27,4).Value = '''3.55'
$newSheet.Cells.Item(27,5).Value = '''68.89'
$newSheet.Cells.Item(27,6).Value = '''2.25'
$newSheet.Cells.Item(27,7).Value = '''0.0799'
$newSheet.Cells.Item(27,8).Value = 8
$newSheet.Cells.Item(28,1).Value = 26
$newSheet.Cells.Item(28,2).Value = '''011741'
$newSheet.Cells.Item(28,3).Value = '博时成长精选混合C'
$newSheet.Cells.Item(28,4).Value = '''1.39'
$newSheet.Cells.Item(28,5).Value = '''75.84'
$newSheet.Cells.Item(28,6).Value = '''5.38'
$newSheet.Cells.Item(28,7).Value = '''0.0748'
$newSheet.Cells.Item(28,8).Value = 4
$newSheet.Cells.Item(29,1).Value = 27
$newSheet.Cells.Item(29,2).Value = '''160528'
$newSheet.Cells.Item(29,3).Value = '博时研究优选3年封闭运作灵活配置混合C'
$newSheet.Cells.Item(29,4).Value = '''0.69'
$newSheet.Cells.Item(29,5).Value = '''95.51'
$newSheet.Cells.Item(29,6).Value = '''9.37'
$newSheet.Cells.Item(29,7).Value = '''0.0647'
$newSheet.Cells.Item(29,8).Value = 2
$newSheet.Cells.Item(30,1).Value = 28
$newSheet.Cells.Item(30,2).Value = '''008967'
$newSheet.Cells.Item(30,3).Value = '博时成长优选两年封闭运作灵活配置混合C'
$newSheet.Cells.Item(30,4).Value = '''0.63'
$newSheet.Cells.Item(30,5).Value = '''80.79'
$newSheet.Cells.Item(30,6).Value = '''8.96'
$newSheet.Cells.Item(30,7).Value = '''0.0564'
$newSheet.Cells.Item(30,8).Value = 1
$newSheet.Cells.Item(31,1).Value = 29
$newSheet.Cells.Item(31,2).Value = '''003580'
$newSheet.Cells.Item(31,3).Value = '泰康沪港深价值优选灵活配置混合'
$newSheet.Cells.Item(31,4).Value = '''1.80'
$newSheet.Cells.Item(31,5).Value = '''85.70'
$newSheet.Cells.Item(31,6).Value = '''2.55'
$newSheet.Cells.Item(31,7).Value = '''0.0459'
$newSheet.Cells.Item(31,8).Value = 6
$newSheet.Cells.Item(32,1).Value = 30
$newSheet.Cells.Item(32,2).Value = '''012345'
$newSheet.Cells.Item(32,3).Value = '嘉实领先优势混合型证券投资基金C'
$newSheet.Cells.Item(32,4).Value = '''0.82'
$newSheet.Cells.Item(32,5).Value = '''82.99'
$newSheet.Cells.Item(32,6).Value = '''4.70'
$newSheet.Cells.Item(32,7).Value = '''0.0385'
$newSheet.Cells.Item(32,8).Value = 5
$newSheet.Cells.Item(33,1).Value = 31
$newSheet.Cells.Item(33,2).Value = '''009592'
$newSheet.Cells.Item(33,3).Value = '博时研究精选一年持有期灵活配置混合C'
$newSheet.Cells.Item(33,4).Value = '''0.73'
$newSheet.Cells.Item(33,5).Value = '''76.70'
$newSheet.Cells.Item(33,6).Value = '''4.89'
$newSheet.Cells.Item(33,7).Value = '''0.0357'
$newSheet.Cells.Item(33,8).Value = 4
$newSheet.Cells.Item(34,1).Value = 32
$newSheet.Cells.Item(34,2).Value = '''012464'
$newSheet.Cells.Item(34,3).Value = '博时成长优势混合型证券投资基金C'
$newSheet.Cells.Item(34,4).Value = '''0.83'
$newSheet.Cells.Item(34,5).Value = '''75.80'
$newSheet.Cells.Item(34,6).Value = '''3.55'
$newSheet.Cells.Item(34,7).Value = '''0.0295'
$newSheet.Cells.Item(34,8).Value = 9
$newSheet.Cells.Item(35,1).Value = 33
$newSheet.Cells.Item(35,2).Value = '''004316'
$newSheet.Cells.Item(35,3).Value = '前海开源沪港深裕鑫灵活配置混合A'
$newSheet.Cells.Item(35,4).Value = '''0.64'
$newSheet.Cells.Item(35,5).Value = '''90.55'
$newSheet.Cells.Item(35,6).Value = '''3.23'
$newSheet.Cells.Item(35,7).Value = '''0.0207'
$newSheet.Cells.Item(35,8).Value = 1
$newSheet.Cells.Item(36,1).Value = 34
$newSheet.Cells.Item(36,2).Value = '''004317'
$newSheet.Cells.Item(36,3).Value = '前海开源沪港深裕鑫灵活配置混合C'
$newSheet.Cells.Item(36,4).Value = '''0.47'
$newSheet.Cells.Item(36,5).Value = '''90.55'
$newSheet.Cells.Item(36,6).Value = '''3.23'
$newSheet.Cells.Item(36,7).Value = '''0.0152'
$newSheet.Cells.Item(36,8).Value = 1
$newSheet.Cells.Item(37,1).Value = 35
$newSheet.Cells.Item(37,2).Value = '''013124'
$newSheet.Cells.Item(37,3).Value = '汇添富精选核心优势一年持有混合C'
$newSheet.Cells.Item(37,4).Value = '''0.30'
$newSheet.Cells.Item(37,5).Value = '''66.61'
$newSheet.Cells.Item(37,6).Value = '''4.26'
$newSheet.Cells.Item(37,7).Value = '''0.0128'
$newSheet.Cells.Item(37,8).Value = 4

# --- Step 3: insert new 2022-Q1 summary row at the top of the "总计" sheet ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Cells.Item(3,1).Copy()
$totalSheet.Cells.Item(2,1).PasteSpecial(-4122)
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 36
$totalSheet.Cells.Item(2,4).Value = 29.92

Write-Output "Done. Sheets:"
foreach ($s in $wb.Worksheets) { Write-Output $s.Name }
